$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 5
    4  = 2
    7  = -1
    8  = 3
    9  = -1
    12 = -2
    28 = -2
    33 = 0
    37 = 0
    39 = -1
    42 = -1
    46 = -3
    52 = 1
    55 = 1
    56 = 1
    57 = 0
    62 = 3
    67 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
